$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: rewrite the sentence describing the slippers/vest outfit
# -----------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "那一身拖鞋背心的行头让人觉得高手都是要这么低调的。",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "那一身拖鞋背心的行头，突然让人觉得有一种《功夫》里的火云邪神的感觉。",
    2)
Write-Output "replace1: $found1"

# -----------------------------------------------------------------
# Change 2: drop the trailing "了" before the closing quote
# -----------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "“好了，不打扰你们收拾行李了，我先回去歇歇了。”",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "“好了，不打扰你们收拾行李了，我先回去歇歇。”",
    2)
Write-Output "replace2: $found2"

# -----------------------------------------------------------------
# Change 3: append the new dorm mix-up dialogue scene at the very
# end of the document body (after the last existing paragraph).
# -----------------------------------------------------------------
function Add-Paragraph([string]$text) {
    $doc = $word.ActiveDocument
    $count = $doc.Paragraphs.Count
    $lastPara = $doc.Paragraphs.Item($count)
    $lastRange = $lastPara.Range
    $lastRange.Collapse(0)
    $lastRange.InsertParagraphAfter()
    $newCount = $doc.Paragraphs.Count
    $newPara = $doc.Paragraphs.Item($newCount)
    $newRange = $newPara.Range
    $newRange.InsertBefore($text)
}

$newParagraphs = @(
    "“这个宿舍满了啊，我们是不是搞错了。”王子明的母亲显然已经看到四张床铺都已经铺满了东西。",
    "“不会啊，报到证上明明写着109宿舍……”王子明又拿出报到单，仔细的看了几遍，确定是109无误。",
    "“啊，不会吧，难道有人搞错了？”张扬走了过来，看了看王子明手上的报到单。",
    "“确实是109，我看看我的错了没。”",
    "于是大家都开始找出自己的报到单，确认有没有出错。",
    "“我的是109。”",
    "“我的也没错。”",
    "“我的……108……”张一帆尴尬的笑了起来，“原来是我走错了，不好意思，不好意思……”",
    "“哈哈，大枣都吃过了，不能退了啊。”古月轩开玩笑的说。",
    "“哈哈……”",
    "“大枣有的是，大家随便吃，”张一帆说着又去包里捧了一大把枣塞给王子明和他的父母。三人抵挡不过山东人民的热情，欣然接受，边吃边夸这枣甜。"
)

foreach ($p in $newParagraphs) {
    Add-Paragraph($p)
}

Write-Output "final paragraph count: $($d.Paragraphs.Count)"
